$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting from column I into the new column J (rows 3-14) by
# copying the existing I-column formatting, then overwrite with the new
# 2022 data values.
$ws.Range("I3:I14").Copy($ws.Range("J3:J14"))

$ws.Range("J4").Value2  = 2022

$ws.Range("J5").Value2  = 96.4
$ws.Range("J6").Value2  = 96.4
$ws.Range("J7").Value2  = 97.9
$ws.Range("J8").Value2  = 95.3
$ws.Range("J9").Value2  = 93.8
$ws.Range("J10").Value2 = 95.5
$ws.Range("J11").Value2 = 94.4
$ws.Range("J12").Value2 = 95
$ws.Range("J13").Value2 = 98.7
$ws.Range("J14").Value2 = 97.3

# J3 stays blank (same as I3), keep just the copied border/format.

# Update the active selection to match the authored workbook state.
$ws.Range("L10").Select() | Out-Null
